# Add I0 and IF columns (I and J) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: copy style/border/font formatting from H1 (same header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-69: I column ("I0") and J column ("IF") values
$data = @(
    @(2, 4, 5),
    @(3, 8, 8),
    @(4, 6, 6),
    @(5, 8, 8),
    @(6, 7, 8),
    @(7, 7, 7),
    @(8, 7, 7),
    @(9, 8, 8),
    @(10, 10, 10),
    @(11, 7, 7),
    @(12, 9, 9),
    @(13, 6, 6),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 5, 6),
    @(17, 5, 6),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 6, 6),
    @(22, 7, 7),
    @(23, 10, 10),
    @(24, 9, 9),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 9, 9),
    @(30, 4, 5),
    @(31, 7, 8),
    @(32, 6, 6),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 7, 7),
    @(36, 6, 7),
    @(37, 6, 6),
    @(38, 9, 9),
    @(39, 7, 7),
    @(40, 8, 8),
    @(41, 7, 7),
    @(42, 8, 8),
    @(43, 9, 9),
    @(44, 8, 8),
    @(45, 7, 7),
    @(46, 8, 8),
    @(47, 9, 9),
    @(48, 6, 7),
    @(49, 8, 8),
    @(50, 8, 8),
    @(51, 8, 8),
    @(52, 7, 8),
    @(53, 6, 6),
    @(54, 8, 8),
    @(55, 8, 8),
    @(56, 6, 7),
    @(57, 8, 8),
    @(58, 7, 7),
    @(59, 6, 6),
    @(60, 9, 9),
    @(61, 6, 6),
    @(62, 8, 8),
    @(63, 7, 7),
    @(64, 7, 7),
    @(65, 6, 6),
    @(66, 6, 6),
    @(67, 7, 7),
    @(68, 7, 7),
    @(69, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Host "Added I0 and IF columns for rows 2-69"